$d = $word.ActiveDocument
$d.Content.Find.Execute("2019/01/16 10:01:53 - Lost user content zone2", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2019/02/08 16:15:54 - Lost user content zone2", 2)
